# Loan RBI, Variable Instalments
# Insert a new (blank) column into the "Repayment schedule" sheet, immediately
# before the existing "Late" column, shifting "Late" / the unlabeled column /
# "Outstanding" one position to the right. Then make the "Repayment schedule"
# tab the active/selected sheet, with cell R8 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a blank column before column N ("Late"); this pushes N->O, O->P, P->Q.
$ws.Columns("N").Insert()

# The newly inserted column N picks up the width of its left neighbour (column M).
$ws.Columns("N").ColumnWidth = 11 - 5/6

# Make "Repayment schedule" the active sheet/tab, with R8 selected.
$ws.Activate()
$ws.Range("R8").Select() | Out-Null
